$wb = $excel.ActiveWorkbook

$decklist = $wb.Worksheets.Item("Decklist")

# Insert a new "Sideboard" sheet right after "Decklist" (so the order becomes
# Decklist, Sideboard, Have) -- this mirrors Excel's real "insert sheet after"
# behaviour and also makes the freshly-added sheet the active tab.
$sideboard = $wb.Worksheets.Add($null, $decklist)
$sideboard.Name = "Sideboard"

# Re-fetch "Have" by name now that the sheet collection has shifted, so we
# don't accidentally keep a stale reference that now points at "Sideboard".
$have = $wb.Worksheets.Item("Have")

# Populate the new Sideboard sheet with its header + cards.
$sideboard.Range("A1").Value = "Name"
$sideboard.Range("B1").Value = "Qty"
$sideboard.Range("A2").Value = "Abrade"
$sideboard.Range("B2").Value = 1
$sideboard.Range("A3").Value = "Into the Flood Maw"
$sideboard.Range("B3").Value = 1
$sideboard.Range("A4").Value = "Spell Pierce"
$sideboard.Range("B4").Value = 2
$sideboard.Range("A6").Value = "Brazen Borrower"
$sideboard.Range("B6").Value = 2
$sideboard.Range("A5").Value = "Mystical Dispute"
$sideboard.Range("B5").Value = 2

# Add the new "Abrade" row to the existing "Have" sheet.
$have.Range("A6").Value = "Abrade"
$have.Range("B6").Value = 1

# Restore each sheet's own selection.
$decklist.Range("A8").Select() | Out-Null
$have.Range("A7").Select() | Out-Null
$sideboard.Range("C6").Select() | Out-Null

# Make the Sideboard tab the active one when the workbook is opened.
$sideboard.Activate()
